# MP131_Transform.xlsx refresh:
#  - drop the now-unused Sheet2 / Sheet3 tabs
#  - refresh the Sheet1 regression coefficients (B2:D4) with the latest run
#  - leave A1:D4 selected (matches the author's last saved selection)

$wb = $excel.ActiveWorkbook

# --- remove the extra (empty) worksheets -------------------------------
foreach ($name in @("Sheet2", "Sheet3")) {
    try {
        $wb.Worksheets($name).Delete() | Out-Null
    } catch {
        # already gone / never existed - nothing to do
    }
}

$ws = $wb.Worksheets("Sheet1")

# --- updated regression coefficients --------------------------------------
$ws.Range("B2").Value = 0.9751372743782489
$ws.Range("C2").Value = 0.08856360304815197
$ws.Range("D2").Value = -0.20313489196415574

$ws.Range("B3").Value = 0.219649687473582
$ws.Range("C3").Value = -0.5076829475142653
$ws.Range("D3").Value = 0.8330738500253057

$ws.Range("B4").Value = -0.029348098931917588
$ws.Range("C4").Value = -0.8569798790043771
$ws.Range("D4").Value = -0.5145135334184373

# --- restore the selection over the whole data range ----------------------
$ws.Activate() | Out-Null
$ws.Range("A1:D4").Select() | Out-Null
